$wb = $excel.ActiveWorkbook

# Update the "展览" sheet (F2: 164 -> 166, F3: 115 -> 116)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 166
$ws1.Range("F3").Value = 116

# Update the "全部类型" sheet (F2: 164 -> 166, F3: 115 -> 116)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 166
$ws4.Range("F3").Value = 116
